# Auto-generated Excel COM-interop script applying the Garuda_Profits sheet updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3151.55
$ws.Range("I19").Value = 3739.1538
$ws.Range("J19").Value = 2060.2856
$ws.Range("K19").Value = 3739.1538
$ws.Range("L19").Value = 2060.2856
$ws.Range("M19").Value = -3564.1538
$ws.Range("N19").Value = -2410.2856

$ws.Range("H92").Value = 556.6667
$ws.Range("I92").Value = 459.78262
$ws.Range("J92").Value = 875
$ws.Range("K92").Value = 459.78262
$ws.Range("L92").Value = 875
$ws.Range("M92").Value = 788.21738
$ws.Range("N92").Value = -3371

$ws.Range("H113").Value = 2899.7144
$ws.Range("I113").Value = 2550
$ws.Range("J113").Value = 3162
$ws.Range("K113").Value = 2550
$ws.Range("L113").Value = 3162
$ws.Range("M113").Value = 704

$ws.Range("H116").Value = 3673
$ws.Range("I116").Value = 1161
$ws.Range("J116").Value = 9953
$ws.Range("K116").Value = 1161
$ws.Range("L116").Value = 9953
$ws.Range("M116").Value = 2281
$ws.Range("N116").Value = -16837

$ws.Range("H135").Value = 37042612
$ws.Range("I135").Value = 50000692
$ws.Range("J135").Value = 19529.285
$ws.Range("K135").Value = 450006228
$ws.Range("L135").Value = 175763.565
$ws.Range("M135").Value = -450003693
$ws.Range("N135").Value = -180833.565

$ws.Range("H140").Value = 33312.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 33312.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 33312.5
$ws.Range("N140").Value = -43672.5

$ws.Range("H141").Value = 4104
$ws.Range("I141").Value = 1683.6
$ws.Range("J141").Value = 7288.737
$ws.Range("K141").Value = 5050.799999999999
$ws.Range("L141").Value = 21866.211
$ws.Range("M141").Value = 129.2000000000007
$ws.Range("N141").Value = -32226.211

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1054.4667
$ws.Range("I2").Value = 531.8182
$ws.Range("J2").Value = 2491.75
$ws.Range("K2").Value = 531.8182
$ws.Range("L2").Value = 2491.75
$ws.Range("M2").Value = -418.8182

$ws.Range("H32").Value = 28469.156
$ws.Range("I32").Value = 29431.326
$ws.Range("J32").Value = 26625
$ws.Range("K32").Value = 29431.326
$ws.Range("L32").Value = 26625
$ws.Range("M32").Value = -29144.326
$ws.Range("N32").Value = -27199

$ws.Range("H45").Value = 786.25
$ws.Range("I45").Value = 798.5714
$ws.Range("J45").Value = 700
$ws.Range("K45").Value = 798.5714
$ws.Range("L45").Value = 700
$ws.Range("M45").Value = -421.5714
$ws.Range("N45").Value = -1454

$ws.Range("H61").Value = 1569.34
$ws.Range("I61").Value = 1248.85
$ws.Range("J61").Value = 2851.3
$ws.Range("K61").Value = 1248.85
$ws.Range("L61").Value = 2851.3
$ws.Range("M61").Value = -1036.85
$ws.Range("N61").Value = -3275.3

$ws.Range("H97").Value = 756.1795
$ws.Range("I97").Value = 738.3333
$ws.Range("J97").Value = 815.6667
$ws.Range("K97").Value = 738.3333
$ws.Range("L97").Value = 815.6667
$ws.Range("M97").Value = -242.3333
$ws.Range("N97").Value = -1807.6667

$ws.Range("H102").Value = 2666.6667
$ws.Range("I102").Value = 4000
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 4000
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -2378

$ws.Range("H116").Value = 1054.4667
$ws.Range("I116").Value = 531.8182
$ws.Range("J116").Value = 2491.75
$ws.Range("K116").Value = 531.8182
$ws.Range("L116").Value = 2491.75
$ws.Range("M116").Value = 1762.1818

$ws.Range("H132").Value = 6041.357
$ws.Range("I132").Value = 7464.472
$ws.Range("J132").Value = 3479.75
$ws.Range("K132").Value = 22393.416
$ws.Range("L132").Value = 10439.25
$ws.Range("M132").Value = -19863.416

$ws.Range("H136").Value = 1569.34
$ws.Range("I136").Value = 1248.85
$ws.Range("J136").Value = 2851.3
$ws.Range("K136").Value = 3746.55
$ws.Range("L136").Value = 8553.900000000001
$ws.Range("M136").Value = -1196.55
$ws.Range("N136").Value = -13653.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1054.4667
$ws.Range("I3").Value = 531.8182
$ws.Range("J3").Value = 2491.75
$ws.Range("K3").Value = 531.8182
$ws.Range("L3").Value = 2491.75
$ws.Range("M3").Value = -417.8182

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H99").Value = 947.75
$ws.Range("I99").Value = 926.6667
$ws.Range("J99").Value = 1011
$ws.Range("K99").Value = 926.6667
$ws.Range("L99").Value = 1011
$ws.Range("M99").Value = 571.3333
$ws.Range("N99").Value = -4007

$ws.Range("H134").Value = 4663.9585
$ws.Range("I134").Value = 5938.643
$ws.Range("J134").Value = 2879.4
$ws.Range("K134").Value = 17815.929
$ws.Range("L134").Value = 8638.200000000001
$ws.Range("M134").Value = -15280.929
$ws.Range("N134").Value = -13708.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 330000
$ws.Range("I6").Value = 330000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 330000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -329887

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1581.4
$ws.Range("I5").Value = 363.07693
$ws.Range("J5").Value = 3844
$ws.Range("K5").Value = 1089.23079
$ws.Range("L5").Value = 11532
$ws.Range("M5").Value = -977.2307900000001
$ws.Range("N5").Value = -11756

$ws.Range("H131").Value = 2796408.2
$ws.Range("I131").Value = 50315
$ws.Range("J131").Value = 2910828.8
$ws.Range("K131").Value = 150945
$ws.Range("L131").Value = 8732486.399999999
$ws.Range("M131").Value = -145905
$ws.Range("N131").Value = -8742566.399999999

$ws.Range("H132").Value = 2187.4119
$ws.Range("I132").Value = 1022.6667
$ws.Range("J132").Value = 2822.7273
$ws.Range("K132").Value = 9204.0003
$ws.Range("L132").Value = 25404.5457
$ws.Range("M132").Value = -6674.0003

$ws.Range("H135").Value = 1581.4
$ws.Range("I135").Value = 363.07693
$ws.Range("J135").Value = 3844
$ws.Range("K135").Value = 3267.69237
$ws.Range("L135").Value = 34596
$ws.Range("M135").Value = -732.6923700000002
$ws.Range("N135").Value = -39666

$ws.Range("H137").Value = 86838776
$ws.Range("I137").Value = 37049930
$ws.Range("J137").Value = 150853020
$ws.Range("K137").Value = 111149790
$ws.Range("L137").Value = 452559060
$ws.Range("M137").Value = -111144690
$ws.Range("N137").Value = -452569260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 12500
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 12500
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 12500
$ws.Range("N18").Value = -12844

$ws.Range("H22").Value = 727.8570999999999
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 498.33334
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 498.33334
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = -1088.33334

$ws.Range("H27").Value = 727.8570999999999
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 498.33334
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 498.33334
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = -712.33334

$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H132").Value = 7980.0625
$ws.Range("I132").Value = 11946.056
$ws.Range("J132").Value = 2880.9285
$ws.Range("K132").Value = 35838.16800000001
$ws.Range("L132").Value = 8642.7855
$ws.Range("M132").Value = -33308.16800000001
$ws.Range("N132").Value = -13702.7855

$ws.Range("H136").Value = 4917.1577
$ws.Range("I136").Value = 5641.857
$ws.Range("J136").Value = 2888
$ws.Range("K136").Value = 16925.571
$ws.Range("L136").Value = 8664
$ws.Range("M136").Value = -14375.571
$ws.Range("N136").Value = -13764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 19097720
$ws.Range("I107").Value = 7812985
$ws.Range("J107").Value = 41667190
$ws.Range("K107").Value = 23438955
$ws.Range("L107").Value = 125001570
$ws.Range("M107").Value = -23437035
$ws.Range("N107").Value = -125005410

$ws.Range("H132").Value = 2151.8215
$ws.Range("I132").Value = 1855.9474
$ws.Range("J132").Value = 2776.4443
$ws.Range("K132").Value = 5567.8422
$ws.Range("L132").Value = 8329.332900000001
$ws.Range("M132").Value = -3037.8422
$ws.Range("N132").Value = -13389.3329

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
